$d = $word.ActiveDocument

$d.Content.Find.Execute("16×40=640", $true, $false, $false, $false, $false, $true, 1, $false, "61×33=2013", 2) | Out-Null
$d.Content.Find.Execute("20×81=1620", $true, $false, $false, $false, $false, $true, 1, $false, "29×82=2378", 2) | Out-Null
$d.Content.Find.Execute("60×29=1740", $true, $false, $false, $false, $false, $true, 1, $false, "38×59=2242", 2) | Out-Null
$d.Content.Find.Execute("69×61=4209", $true, $false, $false, $false, $false, $true, 1, $false, "18×59=1062", 2) | Out-Null
$d.Content.Find.Execute("94×65=6110", $true, $false, $false, $false, $false, $true, 1, $false, "60×37=2220", 2) | Out-Null
$d.Content.Find.Execute("52×99=5148", $true, $false, $false, $false, $false, $true, 1, $false, "77×94=7238", 2) | Out-Null
$d.Content.Find.Execute("92×56=5152", $true, $false, $false, $false, $false, $true, 1, $false, "28×71=1988", 2) | Out-Null
$d.Content.Find.Execute("38×72=2736", $true, $false, $false, $false, $false, $true, 1, $false, "82×91=7462", 2) | Out-Null
$d.Content.Find.Execute("91×28=2548", $true, $false, $false, $false, $false, $true, 1, $false, "64×33=2112", 2) | Out-Null
$d.Content.Find.Execute("70×34=2380", $true, $false, $false, $false, $false, $true, 1, $false, "33×76=2508", 2) | Out-Null
$d.Content.Find.Execute("58×16=928", $true, $false, $false, $false, $false, $true, 1, $false, "36×54=1944", 2) | Out-Null
$d.Content.Find.Execute("28×36=1008", $true, $false, $false, $false, $false, $true, 1, $false, "31×28=868", 2) | Out-Null
$d.Content.Find.Execute("54×62=3348", $true, $false, $false, $false, $false, $true, 1, $false, "48×16=768", 2) | Out-Null
$d.Content.Find.Execute("46×26=1196", $true, $false, $false, $false, $false, $true, 1, $false, "30×63=1890", 2) | Out-Null
$d.Content.Find.Execute("39×93=3627", $true, $false, $false, $false, $false, $true, 1, $false, "28×38=1064", 2) | Out-Null
$d.Content.Find.Execute("61×37=2257", $true, $false, $false, $false, $false, $true, 1, $false, "95×72=6840", 2) | Out-Null
$d.Content.Find.Execute("30×51=1530", $true, $false, $false, $false, $false, $true, 1, $false, "31×67=2077", 2) | Out-Null
$d.Content.Find.Execute("71×66=4686", $true, $false, $false, $false, $false, $true, 1, $false, "68×95=6460", 2) | Out-Null
$d.Content.Find.Execute("89×19=1691", $true, $false, $false, $false, $false, $true, 1, $false, "19×28=532", 2) | Out-Null
$d.Content.Find.Execute("29×31=899", $true, $false, $false, $false, $false, $true, 1, $false, "41×49=2009", 2) | Out-Null
$d.Content.Find.Execute("72×27=1944", $true, $false, $false, $false, $false, $true, 1, $false, "22×45=990", 2) | Out-Null
$d.Content.Find.Execute("45×46=2070", $true, $false, $false, $false, $false, $true, 1, $false, "16×19=304", 2) | Out-Null
$d.Content.Find.Execute("98×98=9604", $true, $false, $false, $false, $false, $true, 1, $false, "94×56=5264", 2) | Out-Null
$d.Content.Find.Execute("62×55=3410", $true, $false, $false, $false, $false, $true, 1, $false, "65×61=3965", 2) | Out-Null
$d.Content.Find.Execute("80×43=3440", $true, $false, $false, $false, $false, $true, 1, $false, "59×46=2714", 2) | Out-Null
